# Holdout Method Results for TS 5050 - add Precision/F1 Score columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. "Model Results" table (rows 2-10): add column P "F1 Score"
# ---------------------------------------------------------------------------

# Header cell P2 - copy the header style from O2 (bold font, border, fill) then set text
$ws.Range("O2").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P2").Value2 = "F1 Score"

# Data cells P3:P10 - F1 = 2*(Precision*Recall)/(Precision+Recall), Precision=L, Recall=M
# Copy number format from the matching "O" cell in the same row (keeps the shaded /
# un-shaded banding already used by the rest of the table), then set the formula.
$rowsTop = @(3,4,5,6,7,8,9,10)
foreach ($r in $rowsTop) {
    $ws.Range("O$r").Copy() | Out-Null
    $ws.Range("P$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("P$r").Formula = "=2*((M$r*L$r)/(M$r+L$r))"
}

# ---------------------------------------------------------------------------
# 2. "Forecast Results" table (rows 14-28): add columns H "Precision" and I "F1 Score"
# ---------------------------------------------------------------------------

# Header cells H14 / I14 - copy header style from G14 (bold font + border)
$ws.Range("G14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value2 = "Precision"

$ws.Range("G14").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null
$ws.Range("I14").Value2 = "F1 Score"

# Data rows 15-28 : Precision (H) and F1 Score (I), computed from TP/FN/TN/FP.
# Values below are taken verbatim (Precision = TP/(TP+FP)*100, F1 computed by the
# same process used for the rest of the workbook).
$hiValues = @{
    15 = @("3.8961038961038961", "7.7790742901594712")
    16 = @("9.67741935483871",   "19.31861630410722")
    17 = @("11.178247734138973", "22.30729810387966")
    18 = @("3.5340314136125657", "7.0627010907949472")
    19 = @("5.2631578947368416", "10.516036956358446")
    20 = @("8.3900226757369616", "16.750797926522853")
    21 = @("8.5168869309838477", "17.007550766072868")
    22 = @("4.4678055190538766", "8.9279852950830438")
    23 = @("9.1999999999999993", "18.365839538458467")
    24 = @("12.903225806451612", "25.748704901449898")
    25 = @("4.1942604856512142", "8.3831542720995387")
    26 = @("4.6296296296296298", "9.2536899088511557")
    27 = @("6.7796610169491522", "13.544591334283338")
    28 = @("7.1005917159763312", "14.183558891318478")
}

foreach ($r in 15..28) {
    $vals = $hiValues[$r]
    # Copy number format from the matching "G" cell in the same row so the
    # header-row shading (rows 22-28 use the darker band) is preserved.
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("H$r").Value2 = [double]$vals[0]

    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("I$r").Value2 = [double]$vals[1]
}

# ---------------------------------------------------------------------------
# 3. Column widths for the newly populated columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 8.3
$ws.Columns.Item(9).ColumnWidth = 7.45
$ws.Columns.Item(16).ColumnWidth = 7.45

# ---------------------------------------------------------------------------
# 4. Selection / page setup
# ---------------------------------------------------------------------------
$ws.Range("O17").Select() | Out-Null
$ws.PageSetup.Orientation = 1   # xlPortrait

Write-Host "Edit complete"
